$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Input Parameters")
$ws2 = $wb.Worksheets.Item("Output Results")

# ---- Sheet 1: Input Parameters ----
$ws1.Range('I2').Value = 4
$ws1.Range('J2').Value = '[0. 0. 0. 0.]'
$ws1.Range('K2').Value = 'Prop'
$ws1.Range('O2').Value = '[''FC'']'
$ws1.Range('A3').Value = 'Summer Tomato - Drip (Gazipur) SD(1)'
$ws1.Range('I3').Value = 4
$ws1.Range('J3').Value = '[0. 0. 0. 0.]'
$ws1.Range('K3').Value = 'Prop'
$ws1.Range('O3').Value = '[''FC'']'
$ws1.Range('A4').Value = 'Summer Tomato - Drip (Gazipur) SD(2)'
$ws1.Range('I4').Value = 4
$ws1.Range('J4').Value = '[0. 0. 0. 0.]'
$ws1.Range('K4').Value = 'Prop'
$ws1.Range('O4').Value = '[''FC'']'
$ws1.Range('A5').Value = 'Winter Tomato - Drip (Gazipur) (Y1)'
$ws1.Range('I5').Value = 4
$ws1.Range('J5').Value = '[0. 0. 0. 0.]'
$ws1.Range('K5').Value = 'Prop'
$ws1.Range('O5').Value = '[''FC'']'
$ws1.Range('A6').Value = 'Winter Tomato - Furrow (Gazipur) (SD1)'
$ws1.Range('I6').Value = 4
$ws1.Range('J6').Value = '[0. 0. 0. 0.]'
$ws1.Range('K6').Value = 'Prop'
$ws1.Range('O6').Value = '[''FC'']'
$ws1.Range('A7').Value = 'Winter Tomato - Drip (Gazipur) (Y2)'
$ws1.Range('I7').Value = 4
$ws1.Range('J7').Value = '[0. 0. 0. 0.]'
$ws1.Range('K7').Value = 'Prop'
$ws1.Range('O7').Value = '[''FC'']'
$ws1.Range('A8').Value = 'Winter Tomato - Furrow (Gazipur) (SD2)'
$ws1.Range('I8').Value = 4
$ws1.Range('J8').Value = '[0. 0. 0. 0.]'
$ws1.Range('K8').Value = 'Prop'
$ws1.Range('O8').Value = '[''FC'']'
$ws1.Range('I9').Value = 4
$ws1.Range('J9').Value = '[0. 0. 0. 0.]'
$ws1.Range('K9').Value = 'Prop'
$ws1.Range('O9').Value = '[''FC'']'
$ws1.Range('I10').Value = 4
$ws1.Range('J10').Value = '[0. 0. 0. 0.]'
$ws1.Range('K10').Value = 'Prop'
$ws1.Range('O10').Value = '[''FC'']'
$ws1.Range('I11').Value = 4
$ws1.Range('J11').Value = '[0. 0. 0. 0.]'
$ws1.Range('K11').Value = 'Prop'
$ws1.Range('O11').Value = '[''FC'']'
$ws1.Range('I12').Value = 4
$ws1.Range('J12').Value = '[0. 0. 0. 0.]'
$ws1.Range('K12').Value = 'Prop'
$ws1.Range('O12').Value = '[''FC'']'
$ws1.Range('I13').Value = 4
$ws1.Range('J13').Value = '[0. 0. 0. 0.]'
$ws1.Range('K13').Value = 'Prop'
$ws1.Range('O13').Value = '[''FC'']'
$ws1.Range('I14').Value = 4
$ws1.Range('J14').Value = '[0. 0. 0. 0.]'
$ws1.Range('K14').Value = 'Prop'
$ws1.Range('O14').Value = '[''FC'']'
$ws1.Range('I15').Value = 4
$ws1.Range('J15').Value = '[0. 0. 0. 0.]'
$ws1.Range('K15').Value = 'Prop'
$ws1.Range('O15').Value = '[''FC'']'
$ws1.Range('I16').Value = 4
$ws1.Range('J16').Value = '[0. 0. 0. 0.]'
$ws1.Range('K16').Value = 'Prop'
$ws1.Range('O16').Value = '[''FC'']'
$ws1.Range('I17').Value = 4
$ws1.Range('J17').Value = '[0. 0. 0. 0.]'
$ws1.Range('K17').Value = 'Prop'
$ws1.Range('O17').Value = '[''FC'']'

# ---- Sheet 2: Output Results ----
$ws2.Range('F2').Value = 8.254857057132657
$ws2.Range('G2').Value = 9.403343776218572
$ws2.Range('F3').Value = 8.277729612791207
$ws2.Range('G3').Value = 113.0564444647878
$ws2.Range('F4').Value = 8.331522595200999
$ws2.Range('G4').Value = 0
$ws2.Range('F5').Value = 8.348821793627463
$ws2.Range('G5').Value = 23.32236107832408
$ws2.Range('A6').Value = 'Summer Tomato - Drip (Gazipur) SD(1)'
$ws2.Range('F6').Value = 8.254857057132657
$ws2.Range('G6').Value = 9.403343776218572
$ws2.Range('A7').Value = 'Summer Tomato - Drip (Gazipur) SD(1)'
$ws2.Range('F7').Value = 8.277729612791207
$ws2.Range('G7').Value = 113.0564444647878
$ws2.Range('A8').Value = 'Summer Tomato - Drip (Gazipur) SD(1)'
$ws2.Range('F8').Value = 8.331522595200999
$ws2.Range('G8').Value = 0
$ws2.Range('A9').Value = 'Summer Tomato - Drip (Gazipur) SD(1)'
$ws2.Range('F9').Value = 8.348821793627463
$ws2.Range('G9').Value = 23.32236107832408
$ws2.Range('A10').Value = 'Summer Tomato - Drip (Gazipur) SD(2)'
$ws2.Range('F10').Value = 8.254857057132657
$ws2.Range('G10').Value = 9.403343776218572
$ws2.Range('A11').Value = 'Summer Tomato - Drip (Gazipur) SD(2)'
$ws2.Range('F11').Value = 8.277729612791207
$ws2.Range('G11').Value = 113.0564444647878
$ws2.Range('A12').Value = 'Summer Tomato - Drip (Gazipur) SD(2)'
$ws2.Range('F12').Value = 8.331522595200999
$ws2.Range('G12').Value = 0
$ws2.Range('A13').Value = 'Summer Tomato - Drip (Gazipur) SD(2)'
$ws2.Range('F13').Value = 8.348821793627463
$ws2.Range('G13').Value = 23.32236107832408
$ws2.Range('A14').Value = 'Winter Tomato - Drip (Gazipur) (Y1)'
$ws2.Range('F14').Value = 8.152824574119061
$ws2.Range('G14').Value = 623.1157872657654
$ws2.Range('A15').Value = 'Winter Tomato - Drip (Gazipur) (Y1)'
$ws2.Range('F15').Value = 8.22490587834845
$ws2.Range('G15').Value = 559.2473247716565
$ws2.Range('A16').Value = 'Winter Tomato - Drip (Gazipur) (Y1)'
$ws2.Range('F16').Value = 8.26548949529284
$ws2.Range('G16').Value = 554.46743171663
$ws2.Range('A17').Value = 'Winter Tomato - Furrow (Gazipur) (SD1)'
$ws2.Range('F17').Value = 8.152824574119061
$ws2.Range('G17').Value = 623.1157872657654
$ws2.Range('A18').Value = 'Winter Tomato - Furrow (Gazipur) (SD1)'
$ws2.Range('F18').Value = 8.22490587834845
$ws2.Range('G18').Value = 559.2473247716565
$ws2.Range('A19').Value = 'Winter Tomato - Furrow (Gazipur) (SD1)'
$ws2.Range('F19').Value = 8.26548949529284
$ws2.Range('G19').Value = 554.46743171663
$ws2.Range('A20').Value = 'Winter Tomato - Drip (Gazipur) (Y2)'
$ws2.Range('F20').Value = 8.152824574119061
$ws2.Range('G20').Value = 623.1157872657654
$ws2.Range('A21').Value = 'Winter Tomato - Drip (Gazipur) (Y2)'
$ws2.Range('F21').Value = 8.22490587834845
$ws2.Range('G21').Value = 559.2473247716565
$ws2.Range('A22').Value = 'Winter Tomato - Drip (Gazipur) (Y2)'
$ws2.Range('F22').Value = 8.26548949529284
$ws2.Range('G22').Value = 554.46743171663
$ws2.Range('A23').Value = 'Winter Tomato - Furrow (Gazipur) (SD2)'
$ws2.Range('F23').Value = 8.152824574119061
$ws2.Range('G23').Value = 623.1157872657654
$ws2.Range('A24').Value = 'Winter Tomato - Furrow (Gazipur) (SD2)'
$ws2.Range('F24').Value = 8.22490587834845
$ws2.Range('G24').Value = 559.2473247716565
$ws2.Range('A25').Value = 'Winter Tomato - Furrow (Gazipur) (SD2)'
$ws2.Range('F25').Value = 8.26548949529284
$ws2.Range('G25').Value = 554.46743171663
$ws2.Range('F26').Value = 8.116140477194804
$ws2.Range('G26').Value = 493.7507086671172
$ws2.Range('F27').Value = 8.151121539390644
$ws2.Range('G27').Value = 454.1084447394541
$ws2.Range('F28').Value = 8.188330610980953
$ws2.Range('G28').Value = 440.1297618031599
$ws2.Range('F29').Value = 8.116140477194804
$ws2.Range('G29').Value = 493.7507086671172
$ws2.Range('F30').Value = 8.151121539390644
$ws2.Range('G30').Value = 454.1084447394541
$ws2.Range('F31').Value = 8.188330610980953
$ws2.Range('G31').Value = 440.1297618031599
$ws2.Range('F32').Value = 8.149854613466564
$ws2.Range('G32').Value = 521.722267415854
$ws2.Range('F33').Value = 8.18529050471156
$ws2.Range('G33').Value = 489.5130707729779
$ws2.Range('F34').Value = 8.234694634428431
$ws2.Range('G34').Value = 475.1946539421569
$ws2.Range('F35').Value = 8.149854613466564
$ws2.Range('G35').Value = 521.722267415854
$ws2.Range('F36').Value = 8.18529050471156
$ws2.Range('G36').Value = 489.5130707729779
$ws2.Range('F37').Value = 8.234694634428431
$ws2.Range('G37').Value = 475.1946539421569
$ws2.Range('F38').Value = 8.160466078310424
$ws2.Range('G38').Value = 396.6686663362698
$ws2.Range('F39').Value = 8.145782130926303
$ws2.Range('G39').Value = 420.0252606894488
$ws2.Range('F40').Value = 8.159362150667588
$ws2.Range('G40').Value = 439.1301208475674
$ws2.Range('F41').Value = 8.160466078310424
$ws2.Range('G41').Value = 396.6686663362698
$ws2.Range('F42').Value = 8.145782130926303
$ws2.Range('G42').Value = 420.0252606894488
$ws2.Range('F43').Value = 8.159362150667588
$ws2.Range('G43').Value = 439.1301208475674
$ws2.Range('F44').Value = 10.98050787905975
$ws2.Range('G44').Value = 789.0064062777657
$ws2.Range('F45').Value = 11.02363371829305
$ws2.Range('G45').Value = 759.754732053014
$ws2.Range('F46').Value = 11.07063181320648
$ws2.Range('G46').Value = 774.2887670850631
$ws2.Range('F47').Value = 6.685001310430317
$ws2.Range('G47').Value = 635.4148958486646
$ws2.Range('F48').Value = 6.663760289770972
$ws2.Range('G48').Value = 576.7753762708826
$ws2.Range('F49').Value = 6.665636245753939
$ws2.Range('G49').Value = 569.842385873952
$ws2.Range('F50').Value = 6.818304607624293
$ws2.Range('G50').Value = 586.3869258097915
$ws2.Range('F51').Value = 6.691117006220352
$ws2.Range('G51').Value = 566.12726644269
$ws2.Range('F52').Value = 6.643738165280067
$ws2.Range('G52').Value = 525.945710868946
$ws2.Range('F53').Value = 6.671830562414732
$ws2.Range('G53').Value = 506.6232797723603
$ws2.Range('F54').Value = 6.743602641799592
$ws2.Range('G54').Value = 527.0375239137037
